$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.233.66'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.43'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.82'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4727'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2892'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06538'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.59'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07942'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.81'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.866.37'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.153'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6810'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '267.20'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.233.05'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.72'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +8.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007386'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.111.70'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.307'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.17%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.180'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.06'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.213'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.85%  '

$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.954'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.395'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09835'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.358'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.95%  '

$ws.Range("E32").Value = '  -0.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.045'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04702'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7022'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.706'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01865'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.602'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.246'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.96'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.924'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8446'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9990'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4156'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.37'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '951.74'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.148'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.190'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.09'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05650'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.22%  '
